# Generate Report for Handoff
# Moves the localization status from "In Translation" to "Ready for
# handoff" and refreshes the associated timestamps on all three sheets
# (Overview, zh-cn, de-de), then lets Excel re-flow the now-wider status
# columns the same way the real report generator does.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-16 18:34:26"
$wsOverview.Columns("E:F").AutoFit()

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-16 18:34:21"
$wsZhCn.Columns("C:C").AutoFit()

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-16 18:34:26"
$wsDeDe.Columns("C:C").AutoFit()
